# "adloori to davuluri completed" - grading pass on the Generic and
# Customer Class sections: fill in the "Total Points" (column E) scores
# and add a grading comment for the toString() deduction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Generic section (rows 3-6) ---
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# --- Customer Class section (rows 10-14) ---
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = "(-1) for wrong output format in toString method"

# --- move the active selection to F12 (was F37) ---
$ws.Range("F12").Select() | Out-Null
